$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update time-slot labels in column C
$ws.Range("C8").Value = "13:5-13:10"
$ws.Range("C9").Value = "13:10-13:15"

# Update the active cell selection to match the saved workbook state
$ws.Range("C18").Select()
